$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column P into the newly introduced column Q,
# row by row, then fill in the 2020 figures reported in the new column.
$qValues = @{
    3 = 2020
    4 = 0.1
    5 = 0.1
    6 = 0.1
    7 = 0
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0.1
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0.1
    26 = 0.2
    27 = 0.1
    28 = 0.3
    29 = 0.4
    30 = 0.2
    31 = 0.2
    32 = 0.2
    33 = 0.1
    35 = 0
    36 = 0.1
    37 = 0.2
}

for ($r = 3; $r -le 37; $r++) {
    $pCell = $ws.Range("P$r")
    $qCell = $ws.Range("Q$r")
    $pCell.Copy() | Out-Null
    $qCell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    if ($qValues.ContainsKey($r)) {
        $qCell.Value = $qValues[$r]
    }
}

$excel.CutCopyMode = $false

# Match the original author selecting P30 after the edit.
$ws.Range("P30").Select() | Out-Null
